# Applies a row-permutation of the weekly data (Fruta / hortaliza, semanal).
# Rows 2-12 keep their static columns (A,B,C,E-J) but columns
# D,K,L,M,N,O,P,Q,R,S,T get reshuffled between rows, as described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for columns D,K,L,M,N,O,P,Q,R,S,T for rows 2..12 after the edit.
$rows = @{
    2  = @{ D = 44174; K = "Castle Brite"; L = "Primera"; M = 75;  N = 9000;  O = 10000; P = 9467;  Q = "`$/caja 10 kilos";          R = "Región de O'Higgins"; S = 947;  T = 10 }
    3  = @{ D = 44544; K = "Castle Brite"; L = "Segunda"; M = 160; N = 16000; O = 17000; P = 16500; Q = "`$/caja 15 kilos";          R = "Región de O'Higgins"; S = 1100; T = 15 }
    4  = @{ D = 44552; K = "Castle Brite"; L = "Primera"; M = 120; N = 15500; O = 16000; P = 15750; Q = "`$/caja 15 kilos";          R = "Región de O'Higgins"; S = 1050; T = 15 }
    5  = @{ D = 44176; K = "Castle Brite"; L = "Primera"; M = 50;  N = 17000; O = 18000; P = 17400; Q = "`$/caja 18 kilos";          R = "Región de O'Higgins"; S = 967;  T = 18 }
    6  = @{ D = 44168; K = "Castle Brite"; L = "Primera"; M = 30;  N = 16000; O = 17000; P = 16500; Q = "`$/caja 16 kilos granel";   R = "Región de Coquimbo";  S = 1031; T = 16 }
    7  = @{ D = 44165; K = "Castle Brite"; L = "Segunda"; M = 60;  N = 16000; O = 17000; P = 16500; Q = "`$/caja 15 kilos granel";   R = "Provincia de Limarí"; S = 1100; T = 15 }
    8  = @{ D = 44181; K = "Modesto";      L = "Primera"; M = 50;  N = 20000; O = 21000; P = 20500; Q = "`$/caja 18 kilos";          R = "Región de Coquimbo";  S = 1139; T = 18 }
    9  = @{ D = 44551; K = "Castle Brite"; L = "Primera"; M = 120; N = 15500; O = 16000; P = 15750; Q = "`$/caja 15 kilos";          R = "Región de O'Higgins"; S = 1050; T = 15 }
    10 = @{ D = 44537; K = "Castle Brite"; L = "Primera"; M = 60;  N = 21000; O = 21500; P = 21250; Q = "`$/caja 15 kilos";          R = "Región de O'Higgins"; S = 1417; T = 15 }
    11 = @{ D = 44189; K = "Dina";         L = "Primera"; M = 80;  N = 16000; O = 17000; P = 16562; Q = "`$/caja 18 kilos";          R = "Región de O'Higgins"; S = 920;  T = 18 }
    12 = @{ D = 44187; K = "Dina";         L = "Primera"; M = 55;  N = 15000; O = 16000; P = 15455; Q = "`$/caja 15 kilos granel";   R = "Región de O'Higgins"; S = 1030; T = 15 }
}

foreach ($r in $rows.Keys) {
    $data = $rows[$r]
    $ws.Cells.Item($r, 4).Value  = $data.D   # D - Fecha
    $ws.Cells.Item($r, 11).Value = $data.K   # K - Variedad
    $ws.Cells.Item($r, 12).Value = $data.L   # L - Calidad
    $ws.Cells.Item($r, 13).Value = $data.M   # M - Volumen
    $ws.Cells.Item($r, 14).Value = $data.N   # N - Precio mínimo
    $ws.Cells.Item($r, 15).Value = $data.O   # O - Precio máximo
    $ws.Cells.Item($r, 16).Value = $data.P   # P - Precio promedio ponderado
    $ws.Cells.Item($r, 17).Value = $data.Q   # Q - Unidad de comercialización
    $ws.Cells.Item($r, 18).Value = $data.R   # R - Origen
    $ws.Cells.Item($r, 19).Value = $data.S   # S - Precio $/Kg
    $ws.Cells.Item($r, 20).Value = $data.T   # T - Kg / unidad
}
